# Apply "Added the stuff at the bottom of the sheets" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J ("pair_kind" generic tag) for the practice rows 2-5
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# New block starting at row 27: "stim details"
$ws.Range("A27").Value = "stim details"

# Header row for the new block
$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# video rows
$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"
$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"
$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

# audio rows
$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"
$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"
